$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh coin prices (D) and 1h volume deltas (E); a handful of rows also had their
# Coin/Link (B/C) swapped because the underlying ranking reordered two coin pairs.
# D-column values are written with a leading apostrophe (forces text / quote-prefix)
# so Excel does not auto-coerce numeric-looking strings (e.g. "1.00") into real numbers;
# ClearFormats() immediately after each D write drops the quote-prefix style flag again
# so the cell ends up back at the default (unstyled) format, same as its neighbours.

$ws.Range('D2').Value = "'66.372.53"
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +4.15%  '
$ws.Range('D3').Value = "'3.490.99"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +2.35%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'592.70"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +4.11%  '
$ws.Range('D6').Value = "'169.69"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +8.14%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').Value = "'3.488.85"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +2.27%  '
$ws.Range('D9').Value = "'0.575"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.25%  '
$ws.Range('D10').Value = "'7.27"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.59%  '
$ws.Range('E11').Value = '  +4.37%  '
$ws.Range('D12').Value = "'0.435"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +2.07%  '
$ws.Range('D13').Value = "'4.090.48"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +2.08%  '
$ws.Range('D14').Value = "'0.134"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.47%  '
$ws.Range('D15').Value = "'27.93"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +3.32%  '
$ws.Range('D16').Value = "'66.292.68"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +3.76%  '
$ws.Range('E17').Value = '  +1.38%  '
$ws.Range('D18').Value = "'3.496.01"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +3.51%  '
$ws.Range('D19').Value = "'6.27"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +2.59%  '
$ws.Range('D20').Value = "'14.01"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +2.75%  '
$ws.Range('D21').Value = "'387.69"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.66%  '
$ws.Range('D22').Value = "'7.98"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +3.08%  '
$ws.Range('D23').Value = "'73.03"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.70%  '
$ws.Range('D24').Value = "'1.00"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('E25').Value = '  +1.89%  '
$ws.Range('E26').Value = '  +6.15%  '
$ws.Range('E27').Value = '  +5.02%  '
$ws.Range('E28').Value = '  +1.67%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('D30').Value = "'6.37"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +5.45%  '
$ws.Range('E31').Value = '  +5.63%  '
$ws.Range('D32').Value = "'2.06"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +3.83%  '
$ws.Range('D33').Value = "'23.41"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.44%  '
$ws.Range('D34').Value = "'7.39"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +6.16%  '
$ws.Range('D35').Value = "'1.00"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('D36').Value = "'1.54"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.69%  '
$ws.Range('D37').Value = "'160.90"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('D38').Value = "'0.900"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +8.58%  '
$ws.Range('E39').Value = '  +5.25%  '
$ws.Range('D40').Value = "'0.0746"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +3.02%  '
$ws.Range('D41').Value = "'27.16"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +5.35%  '
$ws.Range('D42').Value = "'26.36"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.64%  '
$ws.Range('E43').Value = '  +4.47%  '
$ws.Range('D44').Value = "'4.58"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +4.03%  '
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').Value = "'43.43"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.76%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = "'2.794.25"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.80%  '
$ws.Range('E47').Value = '  +2.86%  '
$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').Value = "'354.74"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +9.00%  '
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').Value = "'2.49"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +7.04%  '
$ws.Range('E50').Value = '  +5.57%  '
$ws.Range('D51').Value = "'32.66"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +8.19%  '
